# Update the TPM-derived NATMI ligand-receptor metrics (columns E:T, rows 2-7)
# to reflect the newly recomputed values ("update scripts wuth new tpm").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 0.3333333333333333, 0.04939733333333333, 0.148192, 0.01161064749150641, 0.01161064749150642, 3, 1, 0.110028, 0.330084, 1, 1, 0.005435089791999999, 0.048915808128, 0.01161064749150641, 0.01161064749150642),
    @(3, 1, 0.307438, 0.9223140000000001, 0.07226208385392767, 0.07226208385392767, 3, 1, 0.110028, 0.330084, 1, 1, 0.03382678826400001, 0.304441094376, 0.07226208385392767, 0.07226208385392767),
    @(3, 1, 0.230218, 0.690654, 0.05411182879371945, 0.05411182879371945, 3, 1, 0.110028, 0.330084, 1, 1, 0.025330426104, 0.227973834936, 0.05411182879371945, 0.05411182879371945),
    @(2, 0.6666666666666666, 0.02954266666666666, 0.088628, 0.006943886754192065, 0.006943886754192066, 3, 1, 0.110028, 0.330084, 1, 1, 0.003250520528, 0.029254684752, 0.006943886754192065, 0.006943886754192066),
    @(3, 1, 3.583446, 10.750338, 0.8422747849583385, 0.8422747849583385, 3, 1, 0.110028, 0.330084, 1, 1, 0.394279396488, 3.548514568392, 0.8422747849583385, 0.8422747849583385),
    @(2, 0.6666666666666666, 0.05444366666666667, 0.163331, 0.01279676814831593, 0.01279676814831593, 3, 1, 0.110028, 0.330084, 1, 1, 0.005990327756000001, 0.053912949804, 0.01279676814831593, 0.01279676814831593)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = 5 + $j   # column E = 5
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
